# "Update case 1 and 2" -- corrected retirement and cancer cost.
#
# The bulk of the published diff is Excel-version / resave noise
# (fileVersion rupBuild, absPath, revisionPtr GUID, default font metrics
# -> row heights / dyDescent / default column width, shared-formula
# re-grouping, calcChain regeneration, etc.) that simply falls out of
# opening the workbook in a newer Excel build and isn't something a user
# action reproduces. The real, intentional edits are:
#
#   1. Case_1!C51  : 603590 -> 0        (retirement-step wage corrected)
#   2. Case_1!I53  : 1800000 -> 720000  (cancer cost corrected)
#   3. Case_1!M56  : new note "<--end"  (end-of-table marker)
#
# Everything else in the sheet (J/K columns for rows 51-61, I54, etc.)
# is a pure formula ripple from those two corrected inputs and is left
# for Excel's own recalculation to produce.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Case_1")

# 1) Corrected retirement wage for step 49 (row 51).
$ws1.Range("C51").Value = 0

# 2) Corrected cancer cost for step 51 (row 53).
$ws1.Range("I53").Value = 720000

# 3) New "<--end" marker cell next to the retirement label column.
$ws1.Range("M56").Value = "<--end"

# Cosmetic: the author's last-saved view (zoom + active cell) moved to
# where they were working (row ~36-57, cell I57) while reviewing the fix.
$ws1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 64
$ws1.Range("I57").Select() | Out-Null
